$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Time" values for A2:A42 (replacing the old =5+A(n-1) formulas
# with literal values stepping by 30 instead of 5).
$a = @(0,30,60,90,120,150,180,210,240,270,300,330,360,390,420,450,480,510,540,570,600,630,660,690,720,750,780,810,840,870,900,930,960,990,1020,1050,1080,1110,1140,1170,1200)

# New "Rabbits" values for B2:B42.
$b = @(46.2,75.2,93.1,87,67.5,55,36.5,14.9,26,28.2,33,65.5,78.400000000000006,94.9,88.4,69.3,55.4,26.2,14,16,30,45.6,66.599999999999994,87,96,86.4,61.8,31.8,14.5,6.2,14.1,34.700000000000003,69.3,86.1,89.6,81.8,71.7,56.7,15.2,13.7,10.8)

for ($i = 0; $i -lt $a.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $a[$i]
    $ws.Cells.Item($row, 2).Value = $b[$i]
}

# Update the sheet's saved selection to the whole data range, with no
# distinguished active cell.
$ws.Range("A1:B42").Select()
